$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, shifting existing rows 63-149 down to 64-150
$ws.Rows(63).Insert()

# Populate the newly inserted row 63 with the new weekly data point.
# (Same Mercado/Producto/Variedad/Calidad/price tier as the row that used to be
#  at 63, but a later date and a different volume.)
$ws.Range("A63").Value = 10
$ws.Range("B63").Value = "Vega Modelo de Temuco"
$ws.Range("C63").Value = "La Araucanía"
$ws.Range("D63").Value = Get-Date -Year 2021 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0
$ws.Range("E63").Value = 9
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100102
$ws.Range("H63").Value = "Cítricos"
$ws.Range("I63").Value = 100102006
$ws.Range("J63").Value = "Pomelo"
$ws.Range("K63").Value = "Start Ruby"
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 35
$ws.Range("N63").Value = 12000
$ws.Range("O63").Value = 12000
$ws.Range("P63").Value = 12000
$ws.Range("Q63").Value = "$/bandeja 15 kilos granel"
$ws.Range("R63").Value = "Región de O'Higgins"
$ws.Range("S63").Value = 800
$ws.Range("T63").Value = 15
